$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario tracker")

# Update status column (B) for scenario rows that finished processing:
# previously tagged ".xml", now tagged "own db" (output pulled from GCAM db)
$rows = @(3,4,5,6,7,8,9,11,12,13,14,15,16,17,18,19)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "own db"
}

# Move the active selection to B14
$ws.Range("B14").Select()
